$wb = $excel.ActiveWorkbook

$wsPrincipal = $wb.Worksheets.Item("principal")
$wsPrincipal.Name = "operacion_titulo"

$wsRelacion = $wb.Worksheets.Item("relacion")
$wsRelacion.Name = "operacion_titulo_rel"

$wsRelacion.Activate()
$wsRelacion.Range("K27").Select()
